# Loan RBI, Variable Instalments
#
# The "Repayment Schedule" sheet gains a new (empty) column between the
# existing "In Advance" (M) and "Late" (N) columns - i.e. a new blank
# column N is inserted, pushing the old "Late" column to O and the old
# "Outstanding" column from P to Q. The sheet selection / active cell and
# the active worksheet tab are also updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before the old "Late" column (N), shifting
# the old N (Late) -> O and old P (Outstanding) -> Q.
$ws.Columns("N").Insert()

# The freshly inserted column inherits the width of the column to its
# left ("In Advance", M) - matches Excel's normal insert-column behaviour.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# The "Repayment Schedule" sheet becomes the active/selected tab, with
# R7 selected.
$ws.Activate()
$ws.Range("R7").Select() | Out-Null
